$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 2002, 1998 and 1994 election results (rows 34-51) ---

# 2002
$ws.Range("A34").Value = 2002
$ws.Range("B34").Value = "SPD"
$ws.Range("C34").Value = 251
$ws.Range("D34").Value = "18.488.668"

$ws.Range("A35").Value = 2002
$ws.Range("B35").Value = "CDU"
$ws.Range("C35").Value = 190
$ws.Range("D35").Value = "14.167.561"

$ws.Range("A36").Value = 2002
$ws.Range("B36").Value = "CSU"
$ws.Range("C36").Value = 58
$ws.Range("D36").Value = "4.315.080"

$ws.Range("A37").Value = 2002
$ws.Range("B37").Value = "gruene"
$ws.Range("C37").Value = 55
$ws.Range("D37").Value = "4.110.355"

$ws.Range("A38").Value = 2002
$ws.Range("B38").Value = "FDP"
$ws.Range("C38").Value = 47
$ws.Range("D38").Value = "3.538.815"

$ws.Range("A39").Value = 2002
$ws.Range("B39").Value = "PDS"
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = "1.916.702"

# 1998
$ws.Range("A40").Value = 1998
$ws.Range("B40").Value = "SPD"
$ws.Range("C40").Value = 298
$ws.Range("D40").Value = "20.181.269"

$ws.Range("A41").Value = 1998
$ws.Range("B41").Value = "CDU"
$ws.Range("C41").Value = 198
$ws.Range("D41").Value = "14.004.908"

$ws.Range("A42").Value = 1998
$ws.Range("B42").Value = "CSU"
$ws.Range("C42").Value = 47
$ws.Range("D42").Value = "3.324.480"

$ws.Range("A43").Value = 1998
$ws.Range("B43").Value = "gruene"
$ws.Range("C43").Value = 47
$ws.Range("D43").Value = "3.301.624"

$ws.Range("A44").Value = 1998
$ws.Range("B44").Value = "FDP"
$ws.Range("C44").Value = 43
$ws.Range("D44").Value = "3.080.955"

$ws.Range("A45").Value = 1998
$ws.Range("B45").Value = "PDS"
$ws.Range("C45").Value = 36
$ws.Range("D45").Value = "2.515.454"

# 1994
$ws.Range("A46").Value = 1994
$ws.Range("B46").Value = "CDU"
$ws.Range("C46").Value = 244
$ws.Range("D46").Value = "16.089.960"

$ws.Range("A47").Value = 1994
$ws.Range("B47").Value = "SPD"
$ws.Range("C47").Value = 252
$ws.Range("D47").Value = "17.140.354"

$ws.Range("A48").Value = 1994
$ws.Range("B48").Value = "FDP"
$ws.Range("C48").Value = 47
$ws.Range("D48").Value = "3.258.407"

$ws.Range("A49").Value = 1994
$ws.Range("B49").Value = "CSU"
$ws.Range("C49").Value = 50
$ws.Range("D49").Value = "3.427.196"

$ws.Range("A50").Value = 1994
$ws.Range("B50").Value = "gruene"
$ws.Range("C50").Value = 49
$ws.Range("D50").Value = "3.424.315"

$ws.Range("A51").Value = 1994
$ws.Range("B51").Value = "PDS"
$ws.Range("C51").Value = 30
$ws.Range("D51").Value = "2.066.176"

# --- View state: zoom + selection as left by the editing session ---
$excel.ActiveWindow.Zoom = 166
$ws.Range("C52").Select()
